$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 28: add a new comment cell in column F, and grow the row to two lines ---
$ws.Range("F28").Value = "seemed like the CCY works fine.. But say some errors on random testing though"
$ws.Rows.Item(28).RowHeight = 28.8

# --- New rows 32-34: beta-testing signoff entries ---
$ws.Range("C32").Value = "NMS priorit for consumption"
$ws.Range("C33").Value = "wait, there seems to be no bugs here.. Lets see what assumptions are broken in her rand test"
$ws.Range("B34").Value = "Signoff"
$ws.Range("B32").Value = "Proper ordering"
$ws.Range("B33").Value = "Comment status message"
$ws.Range("E32").Value = "better PE MS group naming"
$ws.Range("A32").Value = 0.27

$ws.Rows.Item(33).RowHeight = 28.8

# --- Update selection to reflect where the author left off editing ---
[void]$ws.Range("E32").Select()
